# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Reorganizes the worker/period table (rows 16-36) from being grouped by
# worker (descending period) to being grouped by period (ascending),
# with the three workers cycling inside each period group.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Each worker's fixed identity data + Valor Mora / Salario Basico figures.
$carlos   = @("18107361",   "CARLOS EDUARDO PORRAS GIL",          354560, 8864000)
$beatriz  = @("41704570",   "BEATRIZ CRISTINA DIAZGRANADOS RUIZ",  25774,  781242)
$daigoro  = @("1016006162", "DAIGORO ALEXANDER ROA BELTRAN",      150000, 3750000)

$workers = @($carlos, $beatriz, $daigoro)
$periods = @("1607", "1608", "1609", "1610", "1611", "1612", "1701")

$row = 16
foreach ($periodo in $periods) {
    foreach ($worker in $workers) {
        $ws.Cells.Item($row, 3).Value = $worker[0]
        $ws.Cells.Item($row, 4).Value = $worker[1]
        $ws.Cells.Item($row, 5).Value = $periodo
        $ws.Cells.Item($row, 6).Value = $worker[2]
        $ws.Cells.Item($row, 7).Value = $worker[3]
        $row = $row + 1
    }
}
